$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price column (D) cells we touch so that
# numeric-looking values are stored as literal text, matching the
# original inlineStr (text) representation of these cells.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.68'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.354'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05738'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.121'

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8159'

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8700'

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1379'

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03177'

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09405'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.739'

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001534'

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04724'

$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006227'
$ws.Range("E18").Value = '17TigerCashTCH'

$ws.Range("B19").Value = 'BitKan'
$ws.Range("C19").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.001243'
$ws.Range("E19").Value = '18BitKanKAN'

$ws.Range("B20").Value = 'HotbitToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.003870'
$ws.Range("E20").Value = '19HotbitTokenHTB'

$ws.Range("B21").Value = 'NitroEx'
$ws.Range("C21").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.00008797'
$ws.Range("E21").Value = '20NitroExNTX'

$ws.Range("B22").Value = 'LEO'
$ws.Range("C22").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.536'
$ws.Range("E22").Value = '21LEOLEO'

$ws.Range("B23").Value = 'BTSEToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.150'
$ws.Range("E23").Value = '22BTSETokenBTSE'

$ws.Range("B24").Value = 'One'
$ws.Range("C24").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.01012'
$ws.Range("E24").Value = '23OneONEBestin24h'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3174'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1319'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1358'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0003013'
$ws.Range("E28").Value = '27UpBotsUBXT'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03717'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006457'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1058'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002998'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007860'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005259'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.3898'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002254'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002099'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0001999'
